# Update "想去人数" (interested-count) figures that changed between the
# previous and newly generated snapshot of the 苏州-漫展信息 data.
#
# Sheet 1 = "展览"      (Exhibitions)
# Sheet 4 = "全部类型"   (All types) -- contains the same events, just at
#                         different row offsets, so values mirror sheet 1.

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item(1)   # 展览
$wsAllTypes    = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet "展览" (sheet1) ---
$wsExhibitions.Range("F2").Value  = 307
$wsExhibitions.Range("F4").Value  = 16977
$wsExhibitions.Range("F5").Value  = 39
$wsExhibitions.Range("F11").Value = 228
$wsExhibitions.Range("F12").Value = 130
$wsExhibitions.Range("F13").Value = 11748
$wsExhibitions.Range("F15").Value = 17
$wsExhibitions.Range("F16").Value = 1439
$wsExhibitions.Range("F17").Value = 4674
$wsExhibitions.Range("F21").Value = 74

# --- Sheet "全部类型" (sheet4) ---
$wsAllTypes.Range("F2").Value  = 307
$wsAllTypes.Range("F5").Value  = 16977
$wsAllTypes.Range("F6").Value  = 39
$wsAllTypes.Range("F12").Value = 228
$wsAllTypes.Range("F13").Value = 130
$wsAllTypes.Range("F16").Value = 11748
$wsAllTypes.Range("F18").Value = 17
$wsAllTypes.Range("F19").Value = 1439
$wsAllTypes.Range("F20").Value = 4674
$wsAllTypes.Range("F24").Value = 74
